# Master_Test_Template.xlsx - "Feat: Added github workflow" edit
# - flips execution_ind (col S) from Y to N for existing test rows (2-9)
# - appends a new data-validation test case (row 10) for the yellow-cabs
#   trip_data table (FM9 / count_check)
# - re-styles the new source cell (C10) and bumps its row height
# - updates the saved cursor/zoom state to match the authoring session

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_validation")
$ws.Activate()

# --- existing rows: execution_ind Y -> N -----------------------------------
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 19).Value = "N"   # column S
}

# --- new row 10: yellow cabs trip_data count_check --------------------------
# (written in the same left-to-right / dropdown-first order the author used,
# so brand-new shared-string entries land in the same slots)
$ws.Range("A10").Value = "FM9"
$ws.Range("D10").Value = "adls"
$ws.Range("C10").Value = "master/yellowcabs/trip_data"
$ws.Range("F10").Value = "adls"
$ws.Range("E10").Value = "NOT APPL"
$ws.Range("G10").Value = "NOT APPL"
$ws.Range("H10").Value = "NOT APPL"
$ws.Range("I10").Value = "yellow_taxi"
$ws.Range("J10").Value = "table"
$ws.Range("K10").Value = "NOT APPL"
$ws.Range("L10").Value = "snowflake_db"
$ws.Range("M10").Value = "NOT APPL"
$ws.Range("N10").Value = "NOT APPL"
$ws.Range("O10").Value = "id"
$ws.Range("B10").Value = "count_check"

# source cell for the new row gets its own (larger, amber) font
$srcCell = $ws.Range("C10")
$srcCell.Font.Name = "Menlo"
$srcCell.Font.Size = 13
$srcCell.Font.Color = 7901646   # RGB(0x78,0x91,0xCE) == OOXML rgb FFCE9178

$ws.Rows.Item(10).RowHeight = 17

# --- view state: zoom + remembered selection --------------------------------
$excel.ActiveWindow.Zoom = 161
$ws.Range("N9").Select()
